$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "time_taken" column (copy formatting from E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Timestamp values for each data row (stored as text)
$timestamps = @(
    "2021-10-05 10:52:47.044813",
    "2021-10-05 10:52:47.044825",
    "2021-10-05 10:52:47.044829",
    "2021-10-05 10:52:47.044832",
    "2021-10-05 10:52:47.044836",
    "2021-10-05 10:52:47.044839",
    "2021-10-05 10:52:47.044842",
    "2021-10-05 10:52:47.044845",
    "2021-10-05 10:52:47.044848"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
